$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

$ws.Range("B2").Value = 16.150000000000002
$ws.Range("C2").Value = 15.350000000000001
$ws.Range("D2").Value = 16.95
$ws.Range("E2").Value = 13.25

$ws.Range("B3").Value = 11.950000000000001
$ws.Range("C3").Value = 6.8500000000000005
$ws.Range("D3").Value = 19.75
$ws.Range("E3").Value = 23.8

$ws.Range("B1:E3").Select()
